# Apply "started burnt chopped wood and minor edits to aspen meadow" changes
# to this workbook (the "aspen meadow" data workbook).
#
# Sheet layout:
#   1 = "F24 % Cover"
#   2 = "F24 Trees"
#   3 = "Historical Frequency"
#   4 = "Historical Rel Frequency"
#   5 = "Historical Avg % Cover"

$wb = $excel.ActiveWorkbook

$wsCover      = $wb.Worksheets.Item(1)   # F24 % Cover
$wsTrees      = $wb.Worksheets.Item(2)   # F24 Trees
$wsFreq       = $wb.Worksheets.Item(3)   # Historical Frequency
$wsRelFreq    = $wb.Worksheets.Item(4)   # Historical Rel Frequency
$wsAvgCover   = $wb.Worksheets.Item(5)   # Historical Avg % Cover

# --- Add a new "GroundCover" row label in column A, row 2 (the year header
# row) on the three "Historical ..." sheets, and drop the stray orphan "x"
# marker that used to sit in G10 with no header of its own. ---
$wsFreq.Range("A2").Value = "GroundCover"
$wsFreq.Range("G10").ClearContents()

$wsRelFreq.Range("A2").Value = "GroundCover"
$wsRelFreq.Range("G10").ClearContents()

$wsAvgCover.Range("A2").Value = "GroundCover"
$wsAvgCover.Range("G10").ClearContents()

# --- Replace the live formulas in K3:K9 on "Historical Rel Frequency" with
# their rounded, static values (calculation no longer needed / locked in). ---
$wsRelFreq.Range("K3").Value = 10.3
$wsRelFreq.Range("K4").Value = 27.6
$wsRelFreq.Range("K5").Value = 3.4
$wsRelFreq.Range("K6").Value = 3.4
$wsRelFreq.Range("K7").Value = 20.7
$wsRelFreq.Range("K8").Value = 27.6
$wsRelFreq.Range("K9").Value = 6.9

# --- Replay the user's final navigation across the sheets so each sheet's
# last selection / the workbook's active tab end up matching: the user
# ends up on "Historical Avg % Cover". ---
$wsCover.Range("J11").Select()
$wsFreq.Range("H11").Select()
$wsRelFreq.Range("K3").Select()
$wsAvgCover.Range("D8").Select()
